$d = $word.ActiveDocument

$replacements = @(
    @("46×31=", "75×31="),
    @("55×38=", "79×79="),
    @("22×28=", "26×43="),
    @("34×68=", "68×36="),
    @("27×52=", "51×39="),
    @("15×30=", "35×73="),
    @("73×99=", "64×57="),
    @("81×88=", "24×85="),
    @("40×66=", "53×39="),
    @("37×18=", "93×47="),
    @("41×97=", "88×21="),
    @("52×23=", "46×49="),
    @("51×68=", "40×82="),
    @("13×83=", "51×35="),
    @("36×51=", "45×59="),
    @("94×93=", "72×76="),
    @("78×77=", "61×75="),
    @("92×29=", "87×69="),
    @("85×69=", "99×59="),
    @("72×98=", "91×25="),
    @("66×81=", "37×73="),
    @("41×56=", "21×18="),
    @("64×16=", "47×89="),
    @("97×53=", "55×82="),
    @("24×59=", "30×46=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
